$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "runs/balls/fours/sixes" figures (columns C:F) for this player's three
# innings rows (2-4) were corrected. The values are numeric-looking but must
# stay stored as text (as in the original file), so force a text number
# format on each touched cell before writing the new value - this prevents
# Excel's automatic "looks like a number" conversion from turning them into
# real numbers.

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2: runs 11->47, balls 13->39, fours unchanged (1), sixes 0->3
Set-TextValue "C2" "47"
Set-TextValue "D2" "39"
Set-TextValue "F2" "3"

# Row 3: runs 47->1, balls 39->2, fours 1->0, sixes 3->0
Set-TextValue "C3" "1"
Set-TextValue "D3" "2"
Set-TextValue "E3" "0"
Set-TextValue "F3" "0"

# Row 4: runs 1->11, balls 2->13, fours 0->1, sixes unchanged (0)
Set-TextValue "C4" "11"
Set-TextValue "D4" "13"
Set-TextValue "E4" "1"
